# Apply updated odds values per the 2024-10-21 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AT2").Value = 2.63
$ws.Range("M2").Value = 1.07
$ws.Range("O2").Value = 1.36
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.73
$ws.Range("U2").Value = 1.8
$ws.Range("V2").Value = 1.95
$ws.Range("M3").Value = 1.08
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3.2
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.67
$ws.Range("U3").Value = 1.95
$ws.Range("V3").Value = 1.8
$ws.Range("J4").Value = 2.38
$ws.Range("M4").Value = 1.08
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.75
$ws.Range("Q4").Value = 2.3
$ws.Range("R4").Value = 1.62
$ws.Range("V4").Value = 1.62
$ws.Range("R5").Value = 1.47
$ws.Range("V5").Value = 1.69
$ws.Range("R6").Value = 1.58
$ws.Range("V6").Value = 1.69
$ws.Range("M7").Value = 1.03
$ws.Range("O7").Value = 1.19
$ws.Range("U7").Value = 1.67
$ws.Range("M8").Value = 1.03
$ws.Range("O8").Value = 1.25
$ws.Range("U8").Value = 1.8
$ws.Range("V8").Value = 1.91
$ws.Range("AB10").Value = 28
$ws.Range("AC10").Value = 12.5
$ws.Range("AD10").Value = 10
$ws.Range("AF10").Value = 120
$ws.Range("AH10").Value = 26
$ws.Range("AI10").Value = 90
$ws.Range("AJ10").Value = 35
$ws.Range("AK10").Value = 400
$ws.Range("AL10").Value = 175
$ws.Range("AN10").Value = 2.9
$ws.Range("AO10").Value = 4.85
$ws.Range("AP10").Value = 17
$ws.Range("AQ10").Value = 11.5
$ws.Range("AR10").Value = 40
$ws.Range("AT10").Value = 3.2
$ws.Range("AU10").Value = 10
$ws.Range("AV10").Value = 110
$ws.Range("AW10").Value = 13
$ws.Range("AX10").Value = 90
$ws.Range("G10").Value = 1.18
$ws.Range("H10").Value = 5.3
$ws.Range("I10").Value = 13
$ws.Range("J10").Value = 1.57
$ws.Range("K10").Value = 2.55
$ws.Range("L10").Value = 10.25
$ws.Range("M10").Value = 1.02
$ws.Range("N10").Value = 14.4
$ws.Range("O10").Value = 1.14
$ws.Range("P10").Value = 4.2
$ws.Range("Q10").Value = 1.55
$ws.Range("R10").Value = 2.15
$ws.Range("S10").Value = 1.27
$ws.Range("T10").Value = 3.42
$ws.Range("W10").Value = 5.9
$ws.Range("X10").Value = 4.85
$ws.Range("Z10").Value = 5.5
$ws.Range("M11").Value = 1.03
$ws.Range("O11").Value = 1.17
$ws.Range("G12").Value = 2.25
$ws.Range("U12").Value = 1.57
$ws.Range("I13").Value = 1.85
$ws.Range("U13").Value = 1.4
$ws.Range("I14").Value = 1.75
$ws.Range("U14").Value = 1.36
$ws.Range("Q16").Value = 1.57
$ws.Range("Q18").Value = 1.41
$ws.Range("U19").Value = 1.54
$ws.Range("U20").Value = 1.63
$ws.Range("AB21").Value = 34
$ws.Range("AG21").Value = 126
$ws.Range("AR21").Value = 81
$ws.Range("G21").Value = 4.75
$ws.Range("I21").Value = 1.57
$ws.Range("U21").Value = 1.5
$ws.Range("V21").Value = 2.37
$ws.Range("AC22").Value = 23
$ws.Range("AD22").Value = 11
$ws.Range("AJ22").Value = 21
$ws.Range("AL22").Value = 41
$ws.Range("AW22").Value = 9
$ws.Range("G22").Value = 1.33
$ws.Range("H22").Value = 5.5
$ws.Range("N22").Value = 23
$ws.Range("U22").Value = 1.54
$ws.Range("V22").Value = 2.25
$ws.Range("Y22").Value = 9.5
$ws.Range("U23").Value = 1.87
$ws.Range("V23").Value = 1.77
$ws.Range("G25").Value = 1.77
$ws.Range("U25").Value = 1.77
$ws.Range("V25").Value = 1.92
$ws.Range("Q26").Value = 1.77
$ws.Range("U26").Value = 1.58
$ws.Range("G27").Value = 1.82
$ws.Range("M27").Value = 1.08
$ws.Range("N27").Value = 8
$ws.Range("O27").Value = 1.4
$ws.Range("P27").Value = 2.75
$ws.Range("Q27").Value = 2.25
$ws.Range("R27").Value = 1.58
$ws.Range("V27").Value = 1.69
$ws.Range("U28").Value = 1.69
$ws.Range("U30").Value = 1.77
$ws.Range("V30").Value = 1.87
$ws.Range("G31").Value = 2.35
$ws.Range("I31").Value = 2.75
$ws.Range("I32").Value = 2.2
$ws.Range("G33").Value = 1.53
$ws.Range("R34").Value = 1.62
$ws.Range("J35").Value = 2.88
